$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right after "2021-Q4" (i.e. right before the
#    "总计" summary sheet).
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Copy the header-row (B1:H1) and the "A" index-column formatting from the
# 2021-Q4 sheet so the new sheet picks up the same bold/centered/bordered
# style used by the other quarterly sheets instead of creating new styles.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Headers
$newSheet.Cells.Item(1,2).Value2 = "基金代码"
$newSheet.Cells.Item(1,3).Value2 = "基金名称"
$newSheet.Cells.Item(1,4).Value2 = "基金规模"
$newSheet.Cells.Item(1,5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1,6).Value2 = "仓位占比"
$newSheet.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value2 = "仓位排名"

# Data row (one fund holding 000009 for 2022-Q1)
$newSheet.Cells.Item(2,1).Value2 = 0

$c = $newSheet.Cells.Item(2,2)
$c.NumberFormat = "@"
$c.Value2 = "159932"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2,3)
$c.NumberFormat = "@"
$c.Value2 = "大成中证500深市ETF"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value2 = "0.42"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2,5)
$c.NumberFormat = "@"
$c.Value2 = "97.30"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2,6)
$c.NumberFormat = "@"
$c.Value2 = "0.96"
$c.Style = "Normal"

$c = $newSheet.Cells.Item(2,7)
$c.NumberFormat = "@"
$c.Value2 = "0.0040"
$c.Style = "Normal"

$newSheet.Cells.Item(2,8).Value2 = 6

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: prepend a "2022-Q1" row, pushing the
#    existing rows down by one.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Give the new bottom row (row 6, for 2021-Q1) the same formatting as the
# other index cells in column A before shifting values into it.
$tot.Range("A2").Copy()
$tot.Range("A6").PasteSpecial(-4122)

# Shift rows 2-5 down to 3-6 (bottom-up so we never overwrite data we still
# need to read).
$tot.Cells.Item(6,1).Value2 = 4
$tot.Cells.Item(6,2).Value2 = $tot.Cells.Item(5,2).Value2
$tot.Cells.Item(6,3).Value2 = $tot.Cells.Item(5,3).Value2
$tot.Cells.Item(6,4).Value2 = $tot.Cells.Item(5,4).Value2

$tot.Cells.Item(5,1).Value2 = 3
$tot.Cells.Item(5,2).Value2 = $tot.Cells.Item(4,2).Value2
$tot.Cells.Item(5,3).Value2 = $tot.Cells.Item(4,3).Value2
$tot.Cells.Item(5,4).Value2 = $tot.Cells.Item(4,4).Value2

$tot.Cells.Item(4,1).Value2 = 2
$tot.Cells.Item(4,2).Value2 = $tot.Cells.Item(3,2).Value2
$tot.Cells.Item(4,3).Value2 = $tot.Cells.Item(3,3).Value2
$tot.Cells.Item(4,4).Value2 = $tot.Cells.Item(3,4).Value2

$tot.Cells.Item(3,1).Value2 = 1
$tot.Cells.Item(3,2).Value2 = $tot.Cells.Item(2,2).Value2
$tot.Cells.Item(3,3).Value2 = $tot.Cells.Item(2,3).Value2
$tot.Cells.Item(3,4).Value2 = $tot.Cells.Item(2,4).Value2

# New first data row: 2022-Q1 totals.
$tot.Cells.Item(2,1).Value2 = 0
$tot.Cells.Item(2,2).Value2 = "2022-Q1"
$tot.Cells.Item(2,3).Value2 = 1
$tot.Cells.Item(2,4).Value2 = 0
